$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("AK_1.png", 0.351, 0.001, 0.451, 0.467, "Tidak Diketahui", "Salah"),
    @("AK_2.png", 0.6889999999999999, 0.002, 0.398, 1, "Akhlak Kamiswara", "Benar"),
    @("AK_3.png", 0.676, 0.002, 0.528, 1, "Akhlak Kamiswara", "Benar"),
    @("AK_4.png", 0.333, 0.001, 0.43, 0.533, "Akhlak Kamiswara", "Benar"),
    @("AK_5.png", 0.332, 0.001, 0.439, 0.667, "Akhlak Kamiswara", "Benar"),
    @("MIB_1.png", 1.4, 0.005, 0.125, 0.4, "Tidak Diketahui", "Salah"),
    @("MIB_2.png", 0.987, 0.003, 0.311, 0.867, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_3.png", 1.613, 0.005, 0.623, 1, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_4.png", 0.982, 0.003, 0.223, 0.9330000000000001, "Muhammad Iqbal Baqi", "Benar"),
    @("MIB_5.png", 1.396, 0.005, 0.246, 0.867, "Muhammad Iqbal Baqi", "Benar"),
    @("AAH_1.png", 0.6919999999999999, 0.002, 0.265, 1, "Andrea Ayunove Hutami", "Benar"),
    @("AAH_2.png", 0.9409999999999999, 0.003, 0.681, 1, "Andrea Ayunove Hutami", "Benar"),
    @("AAH_3.png", 0.724, 0.002, 0.281, 1, "Andrea Ayunove Hutami", "Benar"),
    @("TI_1.png", 0.8129999999999999, 0.003, 0.488, 0.533, "Toni Ismail", "Benar"),
    @("TI_2.png", 0.782, 0.003, 0.413, 0.667, "Toni Ismail", "Benar"),
    @("TI_3.png", 0.593, 0.002, 0.676, 0.667, "Toni Ismail", "Benar"),
    @("TI_4.png", 0.537, 0.002, 0.383, 0.8, "Toni Ismail", "Benar"),
    @("TI_5.png", 0.914, 0.003, 0.449, 0.667, "Toni Ismail", "Benar"),
    @("RAS_1.png", 0.464, 0.002, 0.445, 0.533, "Ridha Ayu Salsabila", "Benar"),
    @("RAS_2.png", 0.885, 0.003, 0.294, 0.867, "Ridha Ayu Salsabila", "Benar"),
    @("RAS_3.png", 0.359, 0.001, 0.343, 0.467, "Tidak Diketahui", "Salah"),
    @("RAS_4.png", 1.079, 0.004, 0.169, 0.4, "Tidak Diketahui", "Salah"),
    @("RAS_5.png", 0.994, 0.003, 0.362, 1, "Ridha Ayu Salsabila", "Benar"),
    @("RR_1.png", 1.345, 0.004, 0.477, 0.9330000000000001, "Rafiqo Rapitasari", "Benar"),
    @("RR_2.png", 1.248, 0.004, 0.52, 0.9330000000000001, "Rafiqo Rapitasari", "Benar"),
    @("RR_3.png", 0.97, 0.003, 0.106, 0.867, "Rafiqo Rapitasari", "Benar"),
    @("RR_4.png", 1.218, 0.004, 0.526, 0.867, "Rafiqo Rapitasari", "Benar"),
    @("RR_5.png", 1.332, 0.004, 0.541, 1, "Rafiqo Rapitasari", "Benar"),
    @("AR_1.png", 0.677, 0.002, 0.397, 1, "Arizli Romadhon", "Benar"),
    @("GA_1.png", 1.513, 0.005, 0.548, 1, "Gege Ardiyansyah", "Benar"),
    @("GA_2.png", 0.487, 0.002, 0.275, 0.9330000000000001, "Gege Ardiyansyah", "Benar"),
    @("GA_3.png", 0.5620000000000001, 0.002, 0.188, 0.733, "Gege Ardiyansyah", "Benar"),
    @("FY_1.png", 0.954, 0.003, 0.308, 0.733, "Fanny Yusuf", "Benar"),
    @("FY_2.png", 1.716, 0.006, 0.47, 0.6, "Fanny Yusuf", "Benar"),
    @("FY_3.png", 1.367, 0.005, 0.468, 0.6, "Fanny Yusuf", "Benar"),
    @("FY_4.png", 1.13, 0.004, 0.376, 0.6, "Fanny Yusuf", "Benar"),
    @("TO_1.png", 0.624, 0.002, 0.397, 0.8, "Tiara Oktavian", "Benar"),
    @("TO_2.png", 0.95, 0.003, 0.374, 1, "Tiara Oktavian", "Benar"),
    @("TO_3.png", 0.649, 0.002, 0.391, 1, "Tiara Oktavian", "Benar"),
    @("TO_4.png", 6.837, 0.023, 0.152, 0.6, "Tiara Oktavian", "Benar"),
    @("TO_5.png", 6.757, 0.022, 0.104, 0.867, "Tiara Oktavian", "Benar"),
    @("TD_1.png", 3.563, 0.012, 0.043, 0.267, "Tidak Diketahui", "Benar"),
    @("TD_2.png", 4.793, 0.016, 0.062, 0.333, "Tidak Diketahui", "Benar"),
    @("TD_3.png", 0.872, 0.003, 0.193, 0.333, "Tidak Diketahui", "Benar"),
    @("TD_4.png", 0.788, 0.003, 0.105, 0.267, "Tidak Diketahui", "Benar")
)

$nRows = $rows.Count
$nCols = 7
$data = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("A1:G45").Value = $data
